$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.009.93'
$ws.Range('E2').Value = '  +1.10%  '

$ws.Range('D3').Value = '2.468.16'
$ws.Range('E3').Value = '  +1.48%  '

$ws.Range('D4').Value = '''1.00'
$ws.Range('E4').Value = '  +0.00%  '

$ws.Range('E5').Value = '  +1.58%  '

$ws.Range('D6').Value = '''146.98'
$ws.Range('E6').Value = '  +1.52%  '

$ws.Range('E7').Value = '  +0.03%  '

$ws.Range('E8').Value = '  +1.70%  '

$ws.Range('D9').Value = '2.466.40'
$ws.Range('E9').Value = '  +1.43%  '

$ws.Range('D10').Value = '''0.111'
$ws.Range('E10').Value = '  +1.09%  '

$ws.Range('E11').Value = '  +0.79%  '

$ws.Range('B12').Value = 'Toncoin'
$ws.Range('C12').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D12').Value = '''5.24'
$ws.Range('E12').Value = '  -0.48%  '

$ws.Range('B13').Value = 'Cardano'
$ws.Range('C13').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D13').Value = '''0.356'
$ws.Range('E13').Value = '  +1.48%  '

$ws.Range('D14').Value = '''27.03'
$ws.Range('E14').Value = '  +1.23%  '

$ws.Range('E15').Value = '  +1.95%  '

$ws.Range('D17').Value = '62.826.66'
$ws.Range('E17').Value = '  +0.99%  '

$ws.Range('D18').Value = '2.473.38'
$ws.Range('E18').Value = '  +1.48%  '

$ws.Range('D19').Value = '''11.44'
$ws.Range('E19').Value = '  +1.42%  '

$ws.Range('D20').Value = '''7.31'
$ws.Range('E20').Value = '  +6.56%  '

$ws.Range('D21').Value = '''326.15'
$ws.Range('E21').Value = '  +0.34%  '

$ws.Range('E22').Value = '  +0.72%  '

$ws.Range('D23').Value = '''1.95'
$ws.Range('E23').Value = '  +12.23%  '

$ws.Range('E24').Value = '  -0.04%  '

$ws.Range('D25').Value = '''65.92'
$ws.Range('E25').Value = '  -2.17%  '

$ws.Range('D26').Value = '''625.13'
$ws.Range('E26').Value = '  +12.78%  '

$ws.Range('E27').Value = '  +8.31%  '

$ws.Range('D28').Value = '''8.48'
$ws.Range('E28').Value = '  -2.56%  '

$ws.Range('B30').Value = 'Binance-PegBSC-USD'
$ws.Range('C30').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D30').Value = '''1.00'
$ws.Range('E30').Value = '  +0.09%  '

$ws.Range('B31').Value = 'Fetch.AI'
$ws.Range('C31').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D31').Value = '''1.49'
$ws.Range('E31').Value = '  +4.46%  '

$ws.Range('D32').Value = '''8.25'
$ws.Range('E32').Value = '  -0.71%  '

$ws.Range('E33').Value = '  -3.63%  '

$ws.Range('D34').Value = '''1.91'
$ws.Range('E34').Value = '  +1.58%  '

$ws.Range('D35').Value = '''5.11'
$ws.Range('E35').Value = '  +5.86%  '

$ws.Range('E36').Value = '  -2.25%  '

$ws.Range('E37').Value = '  -0.02%  '

$ws.Range('E38').Value = '  +0.24%  '

$ws.Range('D39').Value = '''5.41'
$ws.Range('E39').Value = '  -3.31%  '

$ws.Range('D40').Value = '''18.75'
$ws.Range('E40').Value = '  +0.22%  '

$ws.Range('D41').Value = '''147.25'
$ws.Range('E41').Value = '  -1.95%  '

$ws.Range('D42').Value = '''1.78'
$ws.Range('E42').Value = '  -1.41%  '

$ws.Range('D43').Value = '''2.60'
$ws.Range('E43').Value = '  +12.41%  '

$ws.Range('E44').Value = '  +0.09%  '

$ws.Range('D45').Value = '''147.67'
$ws.Range('E45').Value = '  -0.23%  '

$ws.Range('D46').Value = '''3.72'
$ws.Range('E46').Value = '  +1.31%  '

$ws.Range('B47').Value = 'Hedera'
$ws.Range('C47').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D47').Value = '''0.0540'
$ws.Range('E47').Value = '  +0.89%  '

$ws.Range('B48').Value = 'InjectiveProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D48').Value = '''20.75'
$ws.Range('E48').Value = '  +2.25%  '

$ws.Range('D49').Value = '''0.604'
$ws.Range('E49').Value = '  +1.20%  '

$ws.Range('E50').Value = '  +1.29%  '

$ws.Range('D51').Value = '''0.0921'
$ws.Range('E51').Value = '  -0.48%  '
